$d = $word.ActiveDocument

# --- 1. Title: "Play Legends of the Colosseum Megaways Free | Review"
#            -> "Play Legends of the Colosseum Megaways Free - Slot Game Review"
$d.Content.Find.Execute(
    "Play Legends of the Colosseum Megaways Free | Review", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Play Legends of the Colosseum Megaways Free - Slot Game Review", 2) | Out-Null

# --- 2. Copy the (soon to be removed) "Meta description" bold run text - it has the
#        same "empty run + bold run" paragraph shape we need to recreate further down
#        the document, so we re-use it instead of building a run from scratch.
$metaRange = $d.Range(0, $d.Content.End)
$metaRange.Find.Execute("Meta description", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0) | Out-Null
$metaRange.Copy()

# --- 3. Remove the whole "Meta description: ..." paragraph (2nd paragraph).
$d.Paragraphs.Item(2).Range.Delete()

# --- 4. "What we like" bullet list rewrites.
$d.Content.Find.Execute(
    "Megaways feature provides 86,436 ways to win", $true, $false, $false,
    $false, $false, $true, 1, $false, "86,436 ways to win", 2) | Out-Null

$d.Content.Find.Execute(
    "Stunning graphics that change depending on the game", $true, $false, $false,
    $false, $false, $true, 1, $false, "Cascading reel system", 2) | Out-Null

$d.Content.Find.Execute(
    "Gamble feature allows players to double or lose winnings", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Stunning graphics that capture the atmosphere of ancient Rome", 2) | Out-Null

$d.Content.Find.Execute(
    "Theoretical return to the player is 96%", $true, $false, $false,
    $false, $false, $true, 1, $false, "Exciting bonus features and free spins", 2) | Out-Null

# --- 5. "What we don't like" bullet list rewrites.
$d.Content.Find.Execute(
    "Low maximum win of 10,000 times the bet", $true, $false, $false,
    $false, $false, $true, 1, $false, "Limited maximum bet", 2) | Out-Null

$d.Content.Find.Execute(
    "Free spins can be difficult to trigger", $true, $false, $false,
    $false, $false, $true, 1, $false, "Lack of additional bonus games", 2) | Out-Null

# --- 6. Insert a new bold paragraph "Play Legends of the Colosseum Megaways Free -
#        Slot Game Review" right after the "Lack of additional bonus games" bullet
#        (i.e. right before the image-prompt paragraph). Inserting *after* the last
#        bullet (instead of before the prompt paragraph) means the new paragraph
#        inherits plain/non-italic run formatting, which keeps the leading run a
#        true empty run just like the rest of the document.
$lastBulletRange = $d.Range(0, $d.Content.End)
$lastBulletRange.Find.Execute("Lack of additional bonus games", $true, $false, $false,
                               $false, $false, $true, 1, $false, "", 0) | Out-Null
$lastBulletEnd = $lastBulletRange.End

$afterRange = $d.Range($lastBulletEnd, $lastBulletEnd)
$afterRange.InsertParagraphAfter()

$newParaIndex = ($d.Range(0, $lastBulletEnd).Paragraphs.Count) + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newPara.Style = "Normal"

$pasteRange = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$pasteRange.Paste()

$d.Content.Find.Execute(
    "Meta description", $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Legends of the Colosseum Megaways Free - Slot Game Review", 2) | Out-Null

# --- 7. Replace the image-generation prompt text with the new meta description text.
$d.Content.Find.Execute(
    "Prompt: Create a feature image for Legends of the Colosseum Megaways that is in a cartoon style and features a happy Maya warrior with glasses. The image should show the Maya warrior in the Colosseum, surrounded by gladiators and wild tigers, with the reels of the game in the background. The Maya warrior should be holding a sword and shield, looking confident and ready to win big! The image should be colorful and vibrant, with a touch of humor and adventure to reflect the exciting and entertaining nature of the game.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our unbiased review of Legends of the Colosseum Megaways and play for free. Exciting features and stunning graphics!",
    2) | Out-Null
